$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# An empty AQL expression (e.g. an m:dynamicSheet/m:dynamicCell with no
# value) used to generate a table cell whose paragraph still carried a
# (formatted-but-empty) run: <w:p><w:pPr>...</w:pPr><w:r><w:t/></w:r></w:p>.
# That spurious run is what produced the unwanted empty line. Fix: remove
# the run from every cell whose text content is empty, leaving a bare
# <w:p><w:pPr>...</w:pPr></w:p>.
function Clear-EmptyCellRun($table, $rowIndex, $colIndex) {
    $cell = $table.Cell($rowIndex, $colIndex)
    $cr = $cell.Range
    # A collapsed (zero-length) Range.Delete() is a no-op in this engine,
    # so insert a placeholder character first to make the range
    # non-degenerate, then delete everything up to (but excluding) the
    # trailing paragraph mark. This removes the (empty) run while leaving
    # the paragraph - and its pPr - intact.
    $cr.InsertBefore("X")
    $cell2 = $table.Cell($rowIndex, $colIndex)
    $r = $d.Range($cell2.Range.Start, $cell2.Range.End - 1)
    $r.Delete()
}

for ($rowIndex = 1; $rowIndex -le $t.Rows.Count; $rowIndex++) {
    for ($colIndex = 1; $colIndex -le $t.Columns.Count; $colIndex++) {
        $cellText = $t.Cell($rowIndex, $colIndex).Range.Text
        # Cell.Range.Text always includes the trailing cell-mark
        # character(s) (and, for the last cell in a row, the row-end
        # mark), so an "empty" cell's visible text collapses to "".
        $visible = $cellText -replace "[\x07\x0d]", ""
        if ($visible -eq "") {
            Clear-EmptyCellRun $t $rowIndex $colIndex
        }
    }
}
